$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.125.90"
$ws.Range("E2").Value = "  +5.62%  "
$ws.Range("D3").Value = "1.917.60"
$ws.Range("E3").Value = "  +2.51%  "
$ws.Range("E4").Value = "  -0.80%  "
$ws.Range("D5").Value = "'329.39"
$ws.Range("E5").Value = "  +4.36%  "
$ws.Range("E6").Value = "  -0.76%  "
$ws.Range("D7").Value = "'0.5217"
$ws.Range("E7").Value = "  +2.46%  "
$ws.Range("D8").Value = "'0.4076"
$ws.Range("E8").Value = "  +4.49%  "
$ws.Range("D9").Value = "'0.08512"
$ws.Range("E9").Value = "  +1.88%  "
$ws.Range("D10").Value = "'43.01"
$ws.Range("E10").Value = "  +1.76%  "
$ws.Range("E11").Value = "  +1.58%  "
$ws.Range("D12").Value = "'22.33"
$ws.Range("E12").Value = "  +9.66%  "
$ws.Range("D13").Value = "'6.434"
$ws.Range("E13").Value = "  +3.65%  "
$ws.Range("D14").Value = "1.912.13"
$ws.Range("E14").Value = "  +2.21%  "
$ws.Range("D15").Value = "'7.385"
$ws.Range("E15").Value = "  +1.81%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  -0.79%  "
$ws.Range("D17").Value = "'95.01"
$ws.Range("E17").Value = "  +4.05%  "
$ws.Range("D18").Value = "'0.00001113"
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("D19").Value = "'0.06694"
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("E20").Value = "  +4.25%  "
$ws.Range("D21").Value = "'0.9999"
$ws.Range("E21").Value = "  -0.77%  "
$ws.Range("D22").Value = "'6.011"
$ws.Range("E22").Value = "  +1.78%  "
$ws.Range("D23").Value = "30.127.05"
$ws.Range("E23").Value = "  +5.44%  "
$ws.Range("D24").Value = "'11.33"
$ws.Range("E24").Value = "  +2.23%  "
$ws.Range("D25").Value = "'2.217"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("D26").Value = "2.133.79"
$ws.Range("E26").Value = "  +3.06%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'21.15"
$ws.Range("E27").Value = "  +2.79%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "'160.48"
$ws.Range("E28").Value = "  +2.15%  "
$ws.Range("D29").Value = "'2.421"
$ws.Range("E29").Value = "  +0.37%  "
$ws.Range("D30").Value = "'128.92"
$ws.Range("E30").Value = "  +2.34%  "
$ws.Range("D31").Value = "'1.078"
$ws.Range("E31").Value = "  +3.54%  "
$ws.Range("E32").Value = "  +2.26%  "
$ws.Range("D33").Value = "'6.032"
$ws.Range("E33").Value = "  +5.19%  "
$ws.Range("D34").Value = "'3.635"
$ws.Range("E34").Value = "  +0.46%  "
$ws.Range("E35").Value = "  +1.74%  "
$ws.Range("D36").Value = "'0.06603"
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("D37").Value = "'0.2206"
$ws.Range("E37").Value = "  +2.06%  "
$ws.Range("D38").Value = "'1.228"
$ws.Range("E38").Value = "  +3.96%  "
$ws.Range("D39").Value = "'5.172"
$ws.Range("E39").Value = "  +2.57%  "
$ws.Range("D40").Value = "'8.858"
$ws.Range("E40").Value = "  -0.79%  "
$ws.Range("D41").Value = "'0.6531"
$ws.Range("E41").Value = "  +2.59%  "
$ws.Range("D43").Value = "'1.242"
$ws.Range("E43").Value = "  +0.59%  "
$ws.Range("D44").Value = "'0.6149"
$ws.Range("E44").Value = "  +2.57%  "
$ws.Range("D45").Value = "'13.27"
$ws.Range("E45").Value = "  +1.45%  "
$ws.Range("D46").Value = "'3.744"
$ws.Range("E46").Value = "  +1.65%  "
$ws.Range("D47").Value = "'2.077"
$ws.Range("D48").Value = "'1.244"
$ws.Range("E48").Value = "  +2.55%  "
$ws.Range("D49").Value = "'124.32"
$ws.Range("E49").Value = "  +1.54%  "
$ws.Range("D50").Value = "'1.165"
$ws.Range("E50").Value = "  +4.49%  "
$ws.Range("D51").Value = "'79.76"
$ws.Range("E51").Value = "  +4.87%  "
